$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old test rows (3 and 4) entirely - this also drops the now-unused
# shared strings (Test Engineer, Digital, Organisational Performance, Test User 2/3,
# User 2/3, Janine/Devon emails) from the workbook's string table on save.
$ws.Range("A3:L4").ClearContents()

# Overwrite row 2 with the new onboarding record (Isabella Dagg Court / Policy).
$ws.Range("A2").Value = "Isabella"
$ws.Range("B2").Value = "DaggCourt"
$ws.Range("C2").Value = "Isabella Dagg Court"
$ws.Range("D2").Value = "Policy Advisor"
$ws.Range("E2").Value = "Operational Policy "
$ws.Range("F2").Value = "Policy"
$ws.Range("G2").Value = "Wellington"
$ws.Range("H2").Value = "Level 8,7 Waterloo Quay,Pipitea,Wellington"
$ws.Range("I2").Value = "Rebekah.Hood@hud.govt.nz"
$ws.Range("J2").Value = "Permanent"
$ws.Range("L2").Value = 45670

# Leave the cursor where the author's session ended up.
$ws.Range("L2").Select()
